# Update the 2025 FIFA Club World Cup group-stage results ("Fase de Grupos")
# and correct a wrong date on the first phase.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fase de Grupos")

# Correct a wrong date for the game originally on row 7 (July 23 -> June 23)
$ws.Range("C7").Value = 45831

# Fill in match scores (placar) for the games that have been played.
$scores = @{
    "F5"  = 2;  "H5"  = 0;
    "F6"  = 2;  "H6"  = 1;
    "F11" = 0;  "H11" = 1;
    "F12" = 1;  "H12" = 3;
    "F17" = 2;  "H17" = 1;
    "F18" = 6;  "H18" = 0;
    "F23" = 3;  "H23" = 1;
    "F24" = 0;  "H24" = 1;
    "F29" = 0;  "H29" = 0;
    "F30" = 2;  "H30" = 1;
    "F35" = 4;  "H35" = 2;
    "F36" = 3;  "H36" = 4;
    "F40" = 0;  "H40" = 5;
    "F41" = 6;  "H41" = 0;
    "F42" = 4;  "H42" = 1;
    "F47" = 3;  "H47" = 1;
    "F48" = 0;  "H48" = 0;
}

foreach ($addr in $scores.Keys) {
    $ws.Range($addr).Value = $scores[$addr]
}

# Update the active selection to match the author's last position.
$ws.Activate()
$ws.Range("F7").Select()
